# "Season up to 1/17"
# 1. The game that was scheduled next (GSW, 2024-01-15, home) has been
#    played; its final box score is appended as a new row to the "Games"
#    sheet.
# 2. That same game is removed from the top of the "Next" (upcoming
#    games) sheet, and every remaining row shifts up by one.

$wb = $excel.ActiveWorkbook

$gamesSheet = $wb.Worksheets.Item("Games")
$nextSheet  = $wb.Worksheets.Item("Next")

# ---------------------------------------------------------------------
# 1. Append the completed game to the "Games" sheet (new row 41)
# ---------------------------------------------------------------------
$newRow = $gamesSheet.Cells.Item(40, 1).Row + 1

$gamesSheet.Cells.Item($newRow, 1).Value  = 40        # Game
$gamesSheet.Cells.Item($newRow, 2).Value  = 45306     # Date
$gamesSheet.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"
$gamesSheet.Cells.Item($newRow, 3).Value  = 1         # Streak
$gamesSheet.Cells.Item($newRow, 4).Value  = 116       # Pts
$gamesSheet.Cells.Item($newRow, 5).Value  = 102       # Pace
$gamesSheet.Cells.Item($newRow, 6).Value  = 0.488     # eFG
$gamesSheet.Cells.Item($newRow, 7).Value  = 9.6       # TOV
$gamesSheet.Cells.Item($newRow, 8).Value  = 22.4      # ORB
$gamesSheet.Cells.Item($newRow, 9).Value  = 0.372     # FTR
$gamesSheet.Cells.Item($newRow, 10).Value = 113.7     # ORT
$gamesSheet.Cells.Item($newRow, 11).Value = "GSW"     # OppID
$gamesSheet.Cells.Item($newRow, 12).Value = 107       # OppPts
$gamesSheet.Cells.Item($newRow, 13).Value = 0.544     # OppeFG
$gamesSheet.Cells.Item($newRow, 14).Value = 16.8      # OppTOV
$gamesSheet.Cells.Item($newRow, 15).Value = 18.4      # OppORB
$gamesSheet.Cells.Item($newRow, 16).Value = 0.1       # OppFTR
$gamesSheet.Cells.Item($newRow, 17).Value = 104.9     # OppORT
$gamesSheet.Cells.Item($newRow, 18).Value = 1         # Location
$gamesSheet.Cells.Item($newRow, 19).Value = 1         # Target

# ---------------------------------------------------------------------
# 2. Remove the now-played game from the "Next" sheet; remaining rows
#    shift up automatically and the dimension shrinks accordingly.
# ---------------------------------------------------------------------
$nextSheet.Rows.Item(2).Delete()
